$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs","Col18a1","Kdr","ECs",2,0.6666666666666666,10.92859066666667,32.785772,0.2185558471001832,0.2185558471001832,3,1,162.98837,488.96511,0.9909539753179891,0.9909539753179891,1781.233179157213,16031.09861241492,0.2165787855129171,0.2165787855129171),
  @("ECs","Col18a1","Kdr","FAPs",2,0.6666666666666666,10.92859066666667,32.785772,0.2185558471001832,0.2185558471001832,3,1,0.6513563333333333,1.954069,0.003960185305646138,0.003960185305646138,7.118406745140889,64.06566070626801,0.0008655216541491895,0.0008655216541491896),
  @("ECs","Col18a1","Kdr","sCs",2,0.6666666666666666,10.92859066666667,32.785772,0.2185558471001832,0.2185558471001832,3,1,0.8364996666666666,2.509499,0.005085839376364744,0.005085839376364744,9.14176244980311,82.27586204822801,0.001111539933116864,0.001111539933116864),
  @("FAPs","Col18a1","Kdr","ECs",3,1,21.275312,63.825936,0.4254751576214852,0.4254751576214852,3,1,162.98837,488.96511,0.9909539753179891,0.9909539753179891,3467.62842412144,31208.65581709296,0.4216262988440588,0.4216262988440588),
  @("FAPs","Col18a1","Kdr","FAPs",3,1,21.275312,63.825936,0.4254751576214852,0.4254751576214852,3,1,0.6513563333333333,1.954069,0.003960185305646138,0.003960185305646138,13.85780921484267,124.720282933584,0.00168496046713008,0.00168496046713008),
  @("FAPs","Col18a1","Kdr","sCs",3,1,21.275312,63.825936,0.4254751576214852,0.4254751576214852,3,1,0.8364996666666666,2.509499,0.005085839376364744,0.005085839376364744,17.79679139622933,160.171122566064,0.002163898310296346,0.002163898310296346),
  @("sCs","Col18a1","Kdr","ECs",3,1,17.79975,53.39925,0.3559689952783316,0.3559689952783316,3,1,162.98837,488.96511,0.9909539753179891,0.9909539753179891,2901.1522389075,26110.3701501675,0.3527488909610132,0.3527488909610133),
  @("sCs","Col18a1","Kdr","FAPs",3,1,17.79975,53.39925,0.3559689952783316,0.3559689952783316,3,1,0.6513563333333333,1.954069,0.003960185305646138,0.003960185305646138,11.59397989425,104.34581904825,0.001409703184366868,0.001409703184366869),
  @("sCs","Col18a1","Kdr","sCs",3,1,17.79975,53.39925,0.3559689952783316,0.3559689952783316,3,1,0.8364996666666666,2.509499,0.005085839376364744,0.005085839376364744,14.88948494175,134.00536447575,0.001810401132951534,0.001810401132951535)
)

$rowIndex = 2
foreach ($row in $data) {
    for ($colIndex = 1; $colIndex -le $row.Length; $colIndex++) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $row[$colIndex - 1]
    }
    $rowIndex++
}
